$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 391
$ws.Range("A390:L390").Copy()
$ws.Range("A391:L391").PasteSpecial(-4122)
$ws.Range("M390").Copy()
$ws.Range("M391").PasteSpecial(-4122)
$ws.Range("A391").Value = 45193.39238819444
$ws.Range("B391").Value = "kilucas45@gmail.com"
$ws.Range("C391").Value = "미래융합스쿨"
$ws.Range("D391").Value = 20236604
$ws.Range("E391").Value = "김동건"
$ws.Range("F391").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G391").Value = 0.3
$ws.Range("H391").Value = "5:5"
$ws.Range("I391").Value = "20분의 1"
$ws.Range("J391").Value = "20만호, 69만명"
$ws.Range("K391").Value = "경상"
$ws.Range("L391").Value = "Red"
$ws.Range("M391").Value = "모름/무응답"

# Row 392
$ws.Range("A390:L390").Copy()
$ws.Range("A392:L392").PasteSpecial(-4122)
$ws.Range("M390").Copy()
$ws.Range("N392").PasteSpecial(-4122)
$ws.Range("A392").Value = 45193.42050375
$ws.Range("B392").Value = "sun21cc2@naver.com"
$ws.Range("C392").Value = "소프트웨어학부"
$ws.Range("D392").Value = 20235229
$ws.Range("E392").Value = "이민규"
$ws.Range("F392").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G392").Value = 0.7
$ws.Range("H392").Value = "3:7"
$ws.Range("I392").Value = "20분의 1"
$ws.Range("J392").Value = "20만호, 69만명"
$ws.Range("K392").Value = "충청"
$ws.Range("L392").Value = "Black"
$ws.Range("N392").Value = "모름/무응답"

# Row 393
$ws.Range("A390:L390").Copy()
$ws.Range("A393:L393").PasteSpecial(-4122)
$ws.Range("M390").Copy()
$ws.Range("M393").PasteSpecial(-4122)
$ws.Range("A393").Value = 45193.432001712965
$ws.Range("B393").Value = "rdc9118@naver.com"
$ws.Range("C393").Value = "바이오메디컬"
$ws.Range("D393").Value = 20173608
$ws.Range("E393").Value = "김예찬"
$ws.Range("F393").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G393").Value = 0.3
$ws.Range("H393").Value = "4:6"
$ws.Range("I393").Value = "15분의 1"
$ws.Range("J393").Value = "44만호, 153만명"
$ws.Range("K393").Value = "경상"
$ws.Range("L393").Value = "Red"
$ws.Range("M393").Value = "반대한다."

# Row 394
$ws.Range("A390:L390").Copy()
$ws.Range("A394:L394").PasteSpecial(-4122)
$ws.Range("M390").Copy()
$ws.Range("N394").PasteSpecial(-4122)
$ws.Range("A394").Value = 45193.4376615625
$ws.Range("B394").Value = "aissipar3@naver.com"
$ws.Range("C394").Value = "광고홍보학과"
$ws.Range("D394").Value = 20232609
$ws.Range("E394").Value = "김현채"
$ws.Range("F394").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G394").Value = 0.3
$ws.Range("H394").Value = "3:7"
$ws.Range("I394").Value = "15분의 1"
$ws.Range("J394").Value = "44만호, 153만명"
$ws.Range("K394").Value = "경기"
$ws.Range("L394").Value = "Black"
$ws.Range("N394").Value = "모름/무응답"

# Row 395
$ws.Range("A390:L390").Copy()
$ws.Range("A395:L395").PasteSpecial(-4122)
$ws.Range("M390").Copy()
$ws.Range("N395").PasteSpecial(-4122)
$ws.Range("A395").Value = 45193.45792314815
$ws.Range("B395").Value = "rlaek1116@naver.com"
$ws.Range("C395").Value = "법학과"
$ws.Range("D395").Value = 20232704
$ws.Range("E395").Value = "김다정"
$ws.Range("F395").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G395").Value = 0.7
$ws.Range("H395").Value = "6:4"
$ws.Range("I395").Value = "20분의 1"
$ws.Range("J395").Value = "20만호, 69만명"
$ws.Range("K395").Value = "평안"
$ws.Range("L395").Value = "Black"
$ws.Range("N395").Value = "모름/무응답"

# Row 396
$ws.Range("A390:L390").Copy()
$ws.Range("A396:L396").PasteSpecial(-4122)
$ws.Range("M390").Copy()
$ws.Range("M396").PasteSpecial(-4122)
$ws.Range("A396").Value = 45193.49240679399
$ws.Range("B396").Value = "serf0403@naver.com"
$ws.Range("C396").Value = "바이오메디컬학과"
$ws.Range("D396").Value = 20193646
$ws.Range("E396").Value = "정예선"
$ws.Range("F396").Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Range("G396").Value = 0.3
$ws.Range("H396").Value = "5:5"
$ws.Range("I396").Value = "20분의 1"
$ws.Range("J396").Value = "20만호, 69만명"
$ws.Range("K396").Value = "경상"
$ws.Range("L396").Value = "Red"
$ws.Range("M396").Value = "반대한다."

# Row 397
$ws.Range("A390:L390").Copy()
$ws.Range("A397:L397").PasteSpecial(-4122)
$ws.Range("M390").Copy()
$ws.Range("M397").PasteSpecial(-4122)
$ws.Range("A397").Value = 45193.50534855324
$ws.Range("B397").Value = "harin3040@naver.com"
$ws.Range("C397").Value = "심리학과"
$ws.Range("D397").Value = 20232113
$ws.Range("E397").Value = "김현진"
$ws.Range("F397").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G397").Value = 0.1
$ws.Range("H397").Value = "5:5"
$ws.Range("I397").Value = "10분의 1"
$ws.Range("J397").Value = "20만호, 69만명"
$ws.Range("K397").Value = "충청"
$ws.Range("L397").Value = "Red"
$ws.Range("M397").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 398
$ws.Range("A390:L390").Copy()
$ws.Range("A398:L398").PasteSpecial(-4122)
$ws.Range("M390").Copy()
$ws.Range("M398").PasteSpecial(-4122)
$ws.Range("A398").Value = 45193.508059583335
$ws.Range("B398").Value = "shanesun0923@gmail.com"
$ws.Range("C398").Value = "간호학과"
$ws.Range("D398").Value = 20236253
$ws.Range("E398").Value = "선세인"
$ws.Range("F398").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G398").Value = 0.1
$ws.Range("H398").Value = "7:3"
$ws.Range("I398").Value = "20분의 1"
$ws.Range("J398").Value = "44만호, 153만명"
$ws.Range("K398").Value = "평안"
$ws.Range("L398").Value = "Red"
$ws.Range("M398").Value = "모름/무응답"

# Row 399
$ws.Range("A390:L390").Copy()
$ws.Range("A399:L399").PasteSpecial(-4122)
$ws.Range("M390").Copy()
$ws.Range("N399").PasteSpecial(-4122)
$ws.Range("A399").Value = 45193.52365736111
$ws.Range("B399").Value = "milovany03@gmail.com"
$ws.Range("C399").Value = "사회학과"
$ws.Range("D399").Value = 20202223
$ws.Range("E399").Value = "박진옥"
$ws.Range("F399").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G399").Value = 0.1
$ws.Range("H399").Value = "6:4"
$ws.Range("I399").Value = "20분의 1"
$ws.Range("J399").Value = "20만호, 69만명"
$ws.Range("K399").Value = "충청"
$ws.Range("L399").Value = "Black"
$ws.Range("N399").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 400
$ws.Range("A390:L390").Copy()
$ws.Range("A400:L400").PasteSpecial(-4122)
$ws.Range("M390").Copy()
$ws.Range("N400").PasteSpecial(-4122)
$ws.Range("A400").Value = 45193.525144652776
$ws.Range("B400").Value = "bluessin0312@naver.com"
$ws.Range("C400").Value = "국어국문학과"
$ws.Range("D400").Value = 20221054
$ws.Range("E400").Value = "신혜빈"
$ws.Range("F400").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G400").Value = 0.1
$ws.Range("H400").Value = "6:4"
$ws.Range("I400").Value = "10분의 1"
$ws.Range("J400").Value = "20만호, 69만명"
$ws.Range("K400").Value = "전라"
$ws.Range("L400").Value = "Black"
$ws.Range("N400").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 401
$ws.Range("A390:L390").Copy()
$ws.Range("A401:L401").PasteSpecial(-4122)
$ws.Range("M390").Copy()
$ws.Range("N401").PasteSpecial(-4122)
$ws.Range("A401").Value = 45193.55547229167
$ws.Range("B401").Value = "hyeseongi81@gmail.com"
$ws.Range("C401").Value = "체육학과"
$ws.Range("D401").Value = 20224137
$ws.Range("E401").Value = "이혜성"
$ws.Range("F401").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G401").Value = 0.1
$ws.Range("H401").Value = "6:4"
$ws.Range("I401").Value = "30분의 1"
$ws.Range("J401").Value = "130만호, 5백만명"
$ws.Range("K401").Value = "충청"
$ws.Range("L401").Value = "Black"
$ws.Range("N401").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."
